$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Range("A1").Value = "Total time taken for the ride"
$ws.Range("B1").Value = 0.03788287037037037

# Row 2
$ws.Range("A2").Value = "Actual Ampere-hours (Ah)"
$ws.Range("B2").Value = 33.50824527777778

# Row 3
$ws.Range("A3").Value = "Actual Watt-hours (Wh)"
$ws.Range("B3").Value = 1709.207198701666

# Row 4 (unchanged)
$ws.Range("A4").Value = "Starting SoC (Ah)"
$ws.Range("B4").Value = 37.269

# Row 5 (unchanged)
$ws.Range("A5").Value = "Ending SoC (Ah)"
$ws.Range("B5").Value = 3.638

# Row 6
$ws.Range("A6").Value = "Starting SoC (%)"
$ws.Range("B6").Value = 94

# Row 7
$ws.Range("A7").Value = "Ending SoC (%)"
$ws.Range("B7").Value = 9

# Row 8
$ws.Range("A8").Value = "Total distance covered (km)"
$ws.Range("B8").Value = 37.58748425011527

# Row 9
$ws.Range("A9").Value = "Total energy consumption(WH/KM)"
$ws.Range("B9").Value = 45.47277458976055

# Row 10
$ws.Range("A10").Value = "Total SOC consumed(%)"
$ws.Range("B10").Value = 85

# Row 11
$ws.Range("A11").Value = "Mode"
$ws.Range("B11").Value = "Custom mode`n67.83%`nEco mode`n32.17%"

# Row 12
$ws.Range("A12").Value = "Peak Power(kW)"
$ws.Range("B12").Value = 4552.1056

# Row 13
$ws.Range("A13").Value = "Average Power(kW)"
$ws.Range("B13").Value = -1888.047227777232

# Row 14
$ws.Range("A14").Value = "Total Energy Regenerated(kWh)"
$ws.Range("B14").Value = 2.388511011944444

# Row 15
$ws.Range("A15").Value = "Regenerative Effectiveness(%)"
$ws.Range("B15").Value = 0.1395487847036083

# Row 16
$ws.Range("A16").Value = "Highest Cell Voltage(V)"
$ws.Range("B16").Value = 3.35

# Row 17
$ws.Range("A17").Value = "Lowest Cell Voltage(V)"
$ws.Range("B17").Value = 3.047

# Row 18
$ws.Range("A18").Value = "Difference in Cell Voltage(V)"
$ws.Range("B18").Value = 0.3029999999999999

# Row 19
$ws.Range("A19").Value = "Minimum Temperature(C)"
$ws.Range("B19").Value = 37

# Row 20
$ws.Range("A20").Value = "Maximum Temperature(C)"
$ws.Range("B20").Value = 47

# Row 21
$ws.Range("A21").Value = "Difference in Temperature(C)"
$ws.Range("B21").Value = 10

# Row 22
$ws.Range("A22").Value = "Maximum Fet Temperature-BMS(C)"
$ws.Range("B22").Value = 70

# Row 23
$ws.Range("A23").Value = "Maximum Afe Temperature-BMS(C)"
$ws.Range("B23").Value = 67

# Row 24
$ws.Range("A24").Value = "Maximum PCB Temperature-BMS(C)"
$ws.Range("B24").Value = 65

# Row 25
$ws.Range("A25").Value = "Maximum MCU Temperature(C)"
$ws.Range("B25").Value = 68

# Row 26
$ws.Range("A26").Value = "Maximum Motor Temperature(C)"
$ws.Range("B26").Value = 100

# Row 27
$ws.Range("A27").Value = "Abnormal Motor Temperature Detected(C)"
$ws.Range("B27").Value = 0

# Row 28
$ws.Range("A28").Value = "highest cell temp(C)"
$ws.Range("B28").Value = 47

# Row 29
$ws.Range("A29").Value = "lowest cell temp(C)"
$ws.Range("B29").Value = 37

# Row 30
$ws.Range("A30").Value = "Difference between Highest and Lowest Cell Temperature at 100% SOC(C)"
$ws.Range("B30").Value = 10

# Row 31 -- previously "Maximum BMS Temperature in C" / 70, now replaced with Battery Voltage row
$ws.Range("A31").Value = "Battery Voltage(V)"
$ws.Range("B31").Value = 53

# Row 32
$ws.Range("A32").Value = "Total energy charged(kWh)"
$ws.Range("B32").Value = 1.775936999722223

# Row 33
$ws.Range("A33").Value = "Electricity consumption units(kW)"
$ws.Range("B33").Value = 0.0000001507228332588368

# Row 34
$ws.Range("A34").Value = "Idling time percentage"
$ws.Range("B34").Value = 1.307706921664239

# Row 35
$ws.Range("A35").Value = "Time spent in 0-10 km/h"
$ws.Range("B35").Value = 2.9668798682794

# Row 36
$ws.Range("A36").Value = "Time spent in 10-20 km/h"
$ws.Range("B36").Value = 3.726806408713824

# Row 37
$ws.Range("A37").Value = "Time spent in 20-30 km/h"
$ws.Range("B37").Value = 12.33297447913368

# Row 38
$ws.Range("A38").Value = "Time spent in 30-40 km/h"
$ws.Range("B38").Value = 37.59419922740801

# Row 39
$ws.Range("A39").Value = "Time spent in 40-50 km/h"
$ws.Range("B39").Value = 13.05807105313153

# Row 40
$ws.Range("A40").Value = "Time spent in 50-60 km/h"
$ws.Range("B40").Value = 17.36748780951175

# Row 41
$ws.Range("A41").Value = "Time spent in 60-70 km/h"
$ws.Range("B41").Value = 9.153948451649674

# Row 42
$ws.Range("A42").Value = "Time spent in 70-80 km/h"
$ws.Range("B42").Value = 2.409600405294155

# Row 43 (new row)
$ws.Range("A43").Value = "Time spent in 80-90 km/h"
$ws.Range("B43").Value = 0
